$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.882.44'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '3.116.90'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.66'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.13'
$ws.Range('E6').Value = '  +3.15%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('E10').Value = '  -1.48%  '
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.17'
$ws.Range('E13').Value = '  +2.14%  '
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').Value = '3.632.95'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '66.855.70'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.13'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '3.114.06'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.21'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '475.09'
$ws.Range('E20').Value = '  +1.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.709'
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('E22').Value = '  +5.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '83.82'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.23'
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('E25').Value = '  -3.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.16'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.51'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('D33').Value = '0.0₃0950'
$ws.Range('E33').Value = '  -7.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.976'
$ws.Range('E36').Value = '  -2.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '46.92'
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.06'
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.19'
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.123'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.63'
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('D43').Value = '2.811.63'
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '382.64'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0353'
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.54'
$ws.Range('E46').Value = '  -9.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '135.10'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.83'
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('E51').Value = '  -0.80%  '
